# PF-2235 Edit org profile (Classifications) * adjust algorithm to include classifications
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Matchmaker Classifications SubScore" block (mirrors the existing Keyword/Area-of-Interest blocks) ---
# (Written first so the new shared-strings land in the same order as the target file.)

# Section header (row 49) - copy formatting from the analogous "Matchmaker Keyword SubScore:" header (row 41)
$ws.Range("A41").Copy() | Out-Null
$ws.Range("A49").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(49).RowHeight = $ws.Rows.Item(41).RowHeight
$ws.Range("A49").Value = "Matchmaker Classifications SubScore:"

# Data row (row 51) - copy formatting from an analogous data row (row 43)
$ws.Range("A43:B43").Copy() | Out-Null
$ws.Range("A51:B51").PasteSpecial(-4122) | Out-Null
$ws.Range("A51").Value = "Classification match between project and organization"
$ws.Range("B51").Value = 1

# Result row (row 53) - copy formatting from the analogous result row (row 46)
$ws.Range("A46:B46").Copy() | Out-Null
$ws.Range("A53:B53").PasteSpecial(-4122) | Out-Null
$ws.Range("A53").Value = "Matchmaker Classifications SubScore"
$ws.Range("B53").Formula = "=B51"

$excel.CutCopyMode = $false

# --- Row 16: new "Classifications" sub-score line feeding the Taxonomy weighted-average block ---
$ws.Range("B16").Value = "Classifications"
$ws.Range("C16").Formula = "=B53"
$ws.Range("D16").Formula = "=+C16*`$H`$15"

# --- Row 13: SubScore Count is now a fixed literal (4 sub-scores instead of a COUNT() formula) ---
$ws.Range("H13").Value = 4

# --- Row 17: Weighted SubScore total now also includes row 16 ---
$ws.Range("D17").Formula = "=SUM(D13:D16)"

# --- View state: scroll + selection (best effort) ---
$ws.Range("A22").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D50").Select() | Out-Null

$wb.Save()
